$d = $word.ActiveDocument

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the single-line mailing address into two lines:
#    "130 Baroni Ave., San Jose CA 95136" ->
#       "130 Baroni Ave."
#       "San Jose, CA 95136"   (new paragraph, matching formatting)
$rng = $d.Content
$rng.Find.Execute("130 Baroni Ave., San Jose CA 95136", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "130 Baroni Ave."
$addrPara = $rng.Paragraphs(1)
$addrPara.Range.InsertParagraphAfter()
$cityPara = $addrPara.Next()
$cityPara.Range.Text = "San Jose, CA 95136"

# 3. Remove the empty (NoSpacing) paragraph that follows "Board of Directors"
$rng2 = $d.Content
$rng2.Find.Execute("Board of Directors", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bodPara = $rng2.Paragraphs(1)
$emptyPara = $bodPara.Next()
$emptyPara.Range.Delete()
